# Apply the StructureDefinition-admit-count.xlsx update:
#  - rebrand IBM/Alvearie references to LinuxForHealth
#  - bump version 7.0.0 -> 8.0.0
#  - bump the publication date
#  - clear the stale Constraint(s) value on the root "Extension" row (Elements sheet)

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -----------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/admit-count"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet -------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# The "Fixed Value" of Extension.url mirrors the StructureDefinition's own URL.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/admit-count"

# The root "Extension" row's Constraint(s) cell no longer carries the
# ele-1/ext-1 constraint text (it now lives solely on Extension.extension).
$elements.Range("AI2").Value = ""
